# Update cryptocurrency price (Price) and 1h volume-change (Volume(1h)) figures
# on Sheet1, matching the refreshed data feed used to generate this workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel
# (single-dot decimals). Force them to stay plain text, matching the source data,
# then restore the default "Normal" style so no extra formatting is introduced.
$textCells = @(
    @("D5", "399.62"),
    @("D6", "110.33"),
    @("D7", "0.558"),
    @("D9", "0.622"),
    @("D10", "39.47"),
    @("D11", "0.0945"),
    @("D14", "19.15"),
    @("D15", "8.10"),
    @("D17", "1.04"),
    @("D18", "10.92"),
    @("D20", "3.31"),
    @("D22", "12.93"),
    @("D23", "302.98"),
    @("D24", "75.04"),
    @("D25", "3.20"),
    @("D26", "8.07"),
    @("D27", "28.20"),
    @("D28", "4.37"),
    @("D29", "7.31"),
    @("D32", "0.111"),
    @("D33", "11.01"),
    @("D34", "37.78"),
    @("D35", "0.0483"),
    @("D37", "51.56"),
    @("D38", "3.20"),
    @("D39", "3.55"),
    @("D40", "0.999"),
    @("D41", "17.60"),
    @("D42", "1.94"),
    @("D43", "133.50"),
    @("D44", "3.97"),
    @("D46", "0.281"),
    @("D50", "2.38"),
    @("D51", "1.99")
)
foreach ($pair in $textCells) {
    $cell = $ws.Range($pair[0])
    $cell.NumberFormat = "@"
    $cell.Value = $pair[1]
    $cell.Style = "Normal"
}

# Remaining cells are unambiguous text (multi-dot prices, or percentage strings
# with "%" / sign / padding) and can be assigned directly.
$plainCells = @(
    @("D2", "56.872.46"),
    @("E2", "  +11.04%  "),
    @("D3", "3.266.59"),
    @("E3", "  +6.63%  "),
    @("E4", "  +0.35%  "),
    @("E5", "  +2.33%  "),
    @("E6", "  +8.85%  "),
    @("E7", "  +4.83%  "),
    @("E8", "  +0.06%  "),
    @("E9", "  +6.54%  "),
    @("E10", "  +7.28%  "),
    @("E11", "  +11.56%  "),
    @("E12", "  +2.52%  "),
    @("D13", "3.792.52"),
    @("E13", "  +7.01%  "),
    @("E14", "  +4.63%  "),
    @("E15", "  +5.84%  "),
    @("D16", "3.266.70"),
    @("E16", "  +6.83%  "),
    @("E17", "  +2.82%  "),
    @("E18", "  +3.32%  "),
    @("D19", "56.991.58"),
    @("E19", "  +11.24%  "),
    @("E20", "  +4.95%  "),
    @("E21", "  +9.52%  "),
    @("E22", "  +5.49%  "),
    @("E24", "  +7.67%  "),
    @("E25", "  +1.76%  "),
    @("E26", "  +2.58%  "),
    @("E27", "  +5.20%  "),
    @("E28", "  +4.96%  "),
    @("E29", "  +2.86%  "),
    @("E30", "  +5.84%  "),
    @("E31", "  -0.04%  "),
    @("E32", "  +5.46%  "),
    @("E33", "  +3.44%  "),
    @("E34", "  +5.71%  "),
    @("E35", "  -1.27%  "),
    @("E36", "  +5.70%  "),
    @("E37", "  +3.29%  "),
    @("E38", "  +28.32%  "),
    @("E39", "  +6.53%  "),
    @("E40", "  -0.08%  "),
    @("E41", "  +6.16%  "),
    @("E42", "  +5.94%  "),
    @("E43", "  +4.25%  "),
    @("E44", "  +4.71%  "),
    @("E45", "  +4.34%  "),
    @("E46", "  -3.13%  "),
    @("E47", "  +3.25%  "),
    @("D48", "2.159.28"),
    @("E48", "  +4.48%  "),
    @("E49", "  +2.56%  "),
    @("E50", "  -3.29%  "),
    @("E51", "  +38.22%  ")
)
foreach ($pair in $plainCells) {
    $ws.Range($pair[0]).Value = $pair[1]
}
